$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that currently sits after the
#    "... CSS for Laptop" run (it is being relocated to the very end of the
#    document, see step 3 below).
# ---------------------------------------------------------------------------
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
    # no-op if it is not present / already removed
}

# ---------------------------------------------------------------------------
# 2) Helper: append a brand-new paragraph at the very end of the document
#    body by inserting raw WordprocessingML at the (collapsed) end-of-content
#    range. This lets us fully control run/paragraph-mark formatting instead
#    of inheriting it from whatever paragraph used to be last.
# ---------------------------------------------------------------------------
function Add-Paragraph([string]$InnerXml) {
    $w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $ip = $d.Range($d.Content.End, $d.Content.End)
    $ip.InsertXML("<w:p $w>$InnerXml</w:p>") | Out-Null
}

$pPr = '<w:pPr><w:pStyle w:val="Listeavsnitt"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr><w:rPr><w:lang w:eastAsia="nb-NO"/></w:rPr></w:pPr>'
$rPr = '<w:rPr><w:lang w:eastAsia="nb-NO"/></w:rPr>'

function R([string]$Text, [bool]$Preserve) {
    if ($Preserve) {
        return "<w:r>$rPr<w:t xml:space=`"preserve`">$Text</w:t></w:r>"
    } else {
        return "<w:r>$rPr<w:t>$Text</w:t></w:r>"
    }
}

function SpellRun([string]$Text) {
    return '<w:proofErr w:type="spellStart"/>' + (R $Text $false) + '<w:proofErr w:type="spellEnd"/>'
}

# ---------------------------------------------------------------------------
# 3) Three new "to-do" bullet paragraphs (list numId 8 / style Listeavsnitt).
# ---------------------------------------------------------------------------

# --- Paragraph: "Special Events, background elefant image needs a phone version"
$p1 = $pPr
$p1 += R "Special " $true
$p1 += SpellRun "Events"
$p1 += R ", " $true
$p1 += SpellRun "background"
$p1 += R " elefant" $true
$p1 += R " image" $true
$p1 += R " " $true
$p1 += SpellRun "needs"
$p1 += R " a " $true
$p1 += SpellRun "phone"
$p1 += R " " $true
$p1 += SpellRun "version"
Add-Paragraph $p1

# --- Paragraph: "Ant background in exhibition spaces needs to be made smaller and more mobile friendly"
$p2 = $pPr
$p2 += R "Ant " $true
$p2 += SpellRun "background"
$p2 += R " in " $true
$p2 += SpellRun "exhibition"
$p2 += R " " $true
$p2 += SpellRun "spaces"
$p2 += R " " $true
$p2 += SpellRun "needs"
$p2 += R " to be " $true
$p2 += SpellRun "made"
$p2 += R " " $true
$p2 += SpellRun "smaller"
$p2 += R " and more mobile " $true
$p2 += SpellRun "friendly"
Add-Paragraph $p2

# --- Paragraph: "Might consider making the sliders smaller in size" + moved
#     _GoBack bookmark + trailing line break.
$p3 = $pPr
$p3 += SpellRun "Might"
$p3 += R " " $true
$p3 += SpellRun "consider"
$p3 += R " " $true
$p3 += SpellRun "making"
$p3 += R " " $true
$p3 += SpellRun "the"
$p3 += R " " $true
$p3 += SpellRun "sliders"
$p3 += R " " $true
$p3 += SpellRun "smaller"
$p3 += R " in size" $true
$p3 += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$p3 += "<w:r>$rPr<w:br/></w:r>"
Add-Paragraph $p3
